# Refresh the crypto price/volume table (rows 2-51) to the latest scrape.
# Row 40/41 and 44/45 also swap which coin occupies which rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "58.962.63"
$ws.Range("E2").Value = "  -0.83%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.525.86"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5: BNB
$ws.Range("D5").Value = "'537.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "

# Row 6: Solana
$ws.Range("D6").Value = "'137.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.85%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.09%  "

# Row 8: XRP
$ws.Range("E8").Value = "  +0.33%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.524.93"
$ws.Range("E9").Value = "  +0.20%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -0.14%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -2.32%  "

# Row 12: Toncoin
$ws.Range("D12").Value = "'5.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.03%  "

# Row 13: Cardano
$ws.Range("E13").Value = "  -0.50%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.975.30"
$ws.Range("E14").Value = "  +0.12%  "

# Row 15: Avalanche
$ws.Range("E15").Value = "  -1.30%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "58.969.09"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17: ShibaInu
$ws.Range("E17").Value = "  -1.48%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.529.93"
$ws.Range("E18").Value = "  +0.38%  "

# Row 19: Chainlink
$ws.Range("D19").Value = "'11.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "

# Row 20: Polkadot
$ws.Range("E20").Value = "  -0.26%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'324.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22: Dai
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23: Uniswap
$ws.Range("D23").Value = "'5.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.34%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'65.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.88%  "

# Row 25: Polygon
$ws.Range("D25").Value = "'0.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "

# Row 26: Kaspa
$ws.Range("E26").Value = "  -2.17%  "

# Row 27: Binance-PegBSC-USD
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("E28").Value = "  -3.06%  "

# Row 29: PEPE
$ws.Range("E29").Value = "  -1.02%  "

# Row 30: Aptos
$ws.Range("D30").Value = "'6.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "

# Row 31: PancakeSwap
$ws.Range("E31").Value = "  -1.93%  "

# Row 32: Monero
$ws.Range("D32").Value = "'170.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.94%  "

# Row 33: Fetch.AI
$ws.Range("D33").Value = "'1.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.71%  "

# Row 34: USDe
$ws.Range("E34").Value = "  -0.03%  "

# Row 35: ImmutableX
$ws.Range("D35").Value = "'1.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "

# Row 36: EthereumClassic
$ws.Range("D36").Value = "'18.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "

# Row 37: NEARProtocol
$ws.Range("D37").Value = "'4.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.39%  "

# Row 38: Stacks
$ws.Range("E38").Value = "  -3.17%  "

# Row 39: OKB
$ws.Range("D39").Value = "'36.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.81%  "

# Row 40: SuiNetwork
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").Value = "'0.811"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "

# Row 41: Filecoin
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.53%  "

# Row 42: Bittensor
$ws.Range("D42").Value = "'286.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.90%  "

# Row 43: RenderToken
$ws.Range("E43").Value = "  -3.17%  "

# Row 44: FirstDigitalUSD
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45: Aave
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'131.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.26%  "

# Row 46: Mantle
$ws.Range("D46").Value = "'0.609"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.51%  "

# Row 47: WhiteBITCoin
$ws.Range("E47").Value = "  -0.02%  "

# Row 48: Stellar
$ws.Range("E48").Value = "  -1.39%  "

# Row 49: Hedera
$ws.Range("D49").Value = "'0.0508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.91%  "

# Row 50: VeChain
$ws.Range("E50").Value = "  -1.40%  "

# Row 51: InjectiveProtocol
$ws.Range("D51").Value = "'17.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.00%  "
